# "#5: cash & deposit done"
#
# The 存款 (deposit) sheet (sheet3) was left in an unfinished state: its
# row 1 just duplicated row 2's sample values instead of holding real
# column headers, and it was missing the "bank" / "deposit_type" /
# "currency" columns plus the standard trailing metadata block
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that every other sheet in this workbook already
# has. This finishes the sheet to match that standard shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 存款 (deposit)

# ---------------------------------------------------------------
# Row 1: turn the old "sample value" header row into real column names
# ---------------------------------------------------------------
$ws.Cells.Item(1,2).Value = "bank"
$ws.Cells.Item(1,3).Value = "deposit_type"
$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"

$headerCol  = 7,8,9,10,11,12,13
$headerName = "property_category","category","date","legislator_name","legislator_id","source_file","index"
for ($i = 0; $i -lt $headerCol.Length; $i++) {
    $col  = $headerCol[$i]
    $cell = $ws.Cells.Item(1,$col)
    $cell.Value = $headerName[$i]
    # match the bold/bordered header style used by the rest of row 1
    $ws.Cells.Item(1,2).Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Row 2 (index 44) - 永豐商業銀行建成分行 / 活期存款 (checking deposit)
# ---------------------------------------------------------------
$ws.Cells.Item(2,2).Value = "永豐商業銀行建成分行"
$ws.Cells.Item(2,3).Value = "活期存款"
$ws.Cells.Item(2,4).Value = "新臺幣"
$ws.Cells.Item(2,7).Value  = "deposit"
$ws.Cells.Item(2,8).Value  = "normal"
# quote-prefix so the engine stores this as literal text, not an
# auto-converted date serial, then restore the plain data-row style
$ws.Cells.Item(2,9).Value = "'2012-04-12"
$ws.Cells.Item(2,2).Copy()
$ws.Cells.Item(2,9).PasteSpecial(-4122)
$ws.Cells.Item(2,10).Value = "王育敏"
$ws.Cells.Item(2,11).Value = 1728
$ws.Cells.Item(2,12).Value = "tmp48bc1"
$ws.Cells.Item(2,13).Value = 44
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Row 3 (index 45) - 永豐商業銀行建成分行 / 定期存款 (time deposit)
# ---------------------------------------------------------------
$ws.Cells.Item(3,2).Value = "永豐商業銀行建成分行"
$ws.Cells.Item(3,3).Value = "定期存款"
$ws.Cells.Item(3,4).Value = "新臺幣"
$ws.Cells.Item(3,7).Value  = "deposit"
$ws.Cells.Item(3,8).Value  = "normal"
$ws.Cells.Item(3,9).Value = "'2012-04-12"
$ws.Cells.Item(3,2).Copy()
$ws.Cells.Item(3,9).PasteSpecial(-4122)
$ws.Cells.Item(3,10).Value = "王育敏"
$ws.Cells.Item(3,11).Value = 1728
$ws.Cells.Item(3,12).Value = "tmp48bc1"
$ws.Cells.Item(3,13).Value = 45
$excel.CutCopyMode = $false
